# Applies the "New Inpiut data 1" update: bumps a handful of test-data
# identifiers (corporate / user / mobile numbers, usernames) by one and
# moves the saved cursor/active-sheet position to where the author last
# left off editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Data value changes
# ---------------------------------------------------------------------

# BankAdminData
$wsBankAdmin = $wb.Worksheets.Item("BankAdminData")
$wsBankAdmin.Range("B2").Value = "RamMadhav8"
$wsBankAdmin.Range("B3").Value = "Ram8"
$wsBankAdmin.Range("B4").Value = 9916560548

# CorporateRegistrationCreation
$wsCorpReg = $wb.Worksheets.Item("CorporateRegistrationCreation")
$wsCorpReg.Range("B2").Value = 17028
$wsCorpReg.Range("B4").Value = "DMTraders28"
$wsCorpReg.Range("B5").Value = 9916569247
$wsCorpReg.Range("B6").Value = 17028

# WorkGroup
$wsWorkGroup = $wb.Worksheets.Item("WorkGroup")
$wsWorkGroup.Range("B2").Value = 17028

# User Registration
$wsUserReg = $wb.Worksheets.Item("User Registration")
$wsUserReg.Range("B2").Value = 17028
$wsUserReg.Range("B3").Value = 3328
$wsUserReg.Range("B4").Value = "User28"
$wsUserReg.Range("B5").Value = 991656924825

# UserLimit Creation
$wsUserLimit = $wb.Worksheets.Item("UserLimit Creation")
$wsUserLimit.Range("B2").Value = 17028

# Map WorkGroup
$wsMapWorkGroup = $wb.Worksheets.Item("Map WorkGroup")
$wsMapWorkGroup.Range("B2").Value = 17028
$wsMapWorkGroup.Range("B3").Value = 3328

# Entity CIF Mapping
$wsEntityCif = $wb.Worksheets.Item("Entity CIF Mapping")
$wsEntityCif.Range("B2").Value = 17028

# DElinkCIF
$wsDelinkCif = $wb.Worksheets.Item("DElinkCIF")
$wsDelinkCif.Range("B2").Value = 17028

# ---------------------------------------------------------------------
# 2. Per-sheet saved cursor position (selection) updates
# ---------------------------------------------------------------------

$wsBankAdmin.Range("B18").Select()
$wsCorpReg.Range("B18").Select()
$wsWorkGroup.Range("B8").Select()
$wsMapWorkGroup.Range("B18").Select()
$wsEntityCif.Range("B10").Select()
$wsDelinkCif.Range("B12").Select()

# ---------------------------------------------------------------------
# 3. Final active sheet/selection (this also becomes the workbook's
#    active tab), matching activeTab going from 2 (0-based) to 10 -
#    i.e. the "DelinkAccount" sheet.
# ---------------------------------------------------------------------

$wsDelinkAccount = $wb.Worksheets.Item("DelinkAccount")
$wsDelinkAccount.Range("B2").Select()
